$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 20006
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 20006
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 20006
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -20344

$ws.Range("H28").Value = 466.875
$ws.Range("I28").Value = 499.16666
$ws.Range("J28").Value = 370
$ws.Range("K28").Value = 499.16666
$ws.Range("L28").Value = 370
$ws.Range("M28").Value = -14.16665999999998
$ws.Range("N28").Value = -1340

$ws.Range("H33").Value = 146.03572
$ws.Range("I33").Value = 122.809525
$ws.Range("J33").Value = 215.71428
$ws.Range("K33").Value = 122.809525
$ws.Range("L33").Value = 215.71428
$ws.Range("M33").Value = 106.190475
$ws.Range("N33").Value = -673.71428

$ws.Range("H64").Value = 58537.555
$ws.Range("I64").Value = 2791.625
$ws.Range("J64").Value = 103134.3
$ws.Range("K64").Value = 2791.625
$ws.Range("L64").Value = 103134.3
$ws.Range("M64").Value = -2543.625
$ws.Range("N64").Value = -103630.3

$ws.Range("H67").Value = 58537.555
$ws.Range("I67").Value = 2791.625
$ws.Range("J67").Value = 103134.3
$ws.Range("K67").Value = 2791.625
$ws.Range("L67").Value = 103134.3
$ws.Range("M67").Value = -1933.625
$ws.Range("N67").Value = -104850.3

$ws.Range("H76").Value = 4231.3335
$ws.Range("I76").Value = 4231.3335
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 4231.3335
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3916.3335
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 4231.3335
$ws.Range("I79").Value = 4231.3335
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 4231.3335
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -3139.3335
$ws.Range("N79").ClearContents()

$ws.Range("H129").Value = 785.449
$ws.Range("I129").Value = 282.875
$ws.Range("J129").Value = 883.5122
$ws.Range("K129").Value = 848.625
$ws.Range("L129").Value = 2650.5366
$ws.Range("M129").Value = 4151.375
$ws.Range("N129").Value = -12650.5366

$ws.Range("H132").Value = 2049970.6
$ws.Range("I132").Value = 2084120
$ws.Range("J132").Value = 1006
$ws.Range("K132").Value = 6252360
$ws.Range("L132").Value = 3018
$ws.Range("M132").Value = -6249830
$ws.Range("N132").Value = -8078

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 930.3182
$ws.Range("I45").Value = 925.36584
$ws.Range("K45").Value = 925.36584
$ws.Range("M45").Value = -548.36584

$ws.Range("H97").Value = 1328.5714
$ws.Range("I97").Value = 1076.7727
$ws.Range("J97").Value = 2251.8333
$ws.Range("K97").Value = 1076.7727
$ws.Range("L97").Value = 2251.8333
$ws.Range("M97").Value = -580.7727
$ws.Range("N97").Value = -3243.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 600.2381
$ws.Range("I80").Value = 707
$ws.Range("J80").Value = 503.18182
$ws.Range("K80").Value = 707
$ws.Range("L80").Value = 503.18182
$ws.Range("M80").Value = 291
$ws.Range("N80").Value = -2499.18182

$ws.Range("H83").Value = 600.2381
$ws.Range("I83").Value = 707
$ws.Range("J83").Value = 503.18182
$ws.Range("K83").Value = 3535
$ws.Range("L83").Value = 2515.9091
$ws.Range("M83").Value = 1457
$ws.Range("N83").Value = -12499.9091

$ws.Range("H94").Value = 1772.7273
$ws.Range("I94").Value = 1611.1111
$ws.Range("J94").Value = 2500
$ws.Range("K94").Value = 1611.1111
$ws.Range("L94").Value = 2500
$ws.Range("M94").Value = -1160.1111
$ws.Range("N94").Value = -3402

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 33332.668
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 33332.668
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 33332.668
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -34464.668

$ws.Range("H62").Value = 2266.6667
$ws.Range("I62").Value = 2266.6667
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2266.6667
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1642.6667
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 2266.6667
$ws.Range("I65").Value = 2266.6667
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 11333.3335
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -8213.333500000001
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 63.3
$ws.Range("I12").Value = 109.55556
$ws.Range("J12").Value = 25.454546
$ws.Range("K12").Value = 328.66668
$ws.Range("L12").Value = 76.36363800000001
$ws.Range("M12").Value = -155.66668
$ws.Range("N12").Value = -422.363638

$ws.Range("H92").Value = 200
$ws.Range("J92").Value = 200
$ws.Range("L92").Value = 600
$ws.Range("N92").Value = -3096

$ws.Range("H132").Value = 1322.8889
$ws.Range("I132").Value = 1060.8667
$ws.Range("J132").Value = 1650.4166
$ws.Range("K132").Value = 9547.800300000001
$ws.Range("L132").Value = 14853.7494
$ws.Range("M132").Value = -7017.800300000001
$ws.Range("N132").Value = -19913.7494

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2379.32
$ws.Range("I132").Value = 1845.8667
$ws.Range("J132").Value = 3179.5
$ws.Range("K132").Value = 5537.6001
$ws.Range("L132").Value = 9538.5
$ws.Range("M132").Value = -3007.6001
$ws.Range("N132").Value = -14598.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2404.2444
$ws.Range("I68").Value = 937.7692
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 937.7692
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -188.7692
$ws.Range("N68").Value = -4498

$ws.Range("H71").Value = 2404.2444
$ws.Range("I71").Value = 937.7692
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 4688.846
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -944.8459999999995
$ws.Range("N71").Value = -22488

$ws.Range("H100").Value = 2358.889
$ws.Range("I100").Value = 1781.6666
$ws.Range("J100").Value = 2523.8096
$ws.Range("K100").Value = 1781.6666
$ws.Range("L100").Value = 2523.8096
$ws.Range("M100").Value = -1240.6666
$ws.Range("N100").Value = -3605.8096

$ws.Range("H132").Value = 4684.7144
$ws.Range("I132").Value = 2247.9092
$ws.Range("J132").Value = 13619.667
$ws.Range("K132").Value = 6743.7276
$ws.Range("L132").Value = 40859.001
$ws.Range("M132").Value = -4213.7276
$ws.Range("N132").Value = -45919.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1302.909
$ws.Range("I100").Value = 1424.2
$ws.Range("J100").Value = 90
$ws.Range("K100").Value = 2848.4
$ws.Range("L100").Value = 180
$ws.Range("M100").Value = -2307.4
$ws.Range("N100").Value = -1262
